$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for the Tab Name, shifting existing data to B:E
$ws.Columns.Item(1).Insert()

# --- Here-strings holding the long Cypher query texts (verbatim, preserve newlines) ---
$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Border Collie']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Border Collie']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Border Collie']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Border Collie']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@

# --- Row 1: header ---
$ws.Cells.Item(1,1).Value = "TabName"
$ws.Cells.Item(1,2).Value = "query"
$ws.Cells.Item(1,3).Value = "StatQuery"
$ws.Cells.Item(1,4).Value = "dbExcel"
$ws.Cells.Item(1,5).Value = "WebExcel"

# --- Row 2: CasesTab ---
$ws.Cells.Item(2,1).Value = "CasesTab"
$ws.Cells.Item(2,2).Value = $casesQuery
$ws.Cells.Item(2,3).Value = $statQuery
$ws.Cells.Item(2,4).Value = "TC09_Canine_Filter_Breed-BorderCol_Neo4jData.xlsx"
$ws.Cells.Item(2,5).Value = "TC09_Canine_Filter_Breed-BorderCol_WebData.xlsx"
$ws.Rows.Item(2).RowHeight = 275.5

# --- Row 3: SamplesTab ---
$ws.Cells.Item(3,1).Value = "SamplesTab"
$ws.Cells.Item(3,2).Value = $samplesQuery
$ws.Cells.Item(3,3).Value = $statQuery
$ws.Cells.Item(3,4).Value = "TC09_Canine_Filter_Breed-BorderCol_Neo4jData.xlsx"
$ws.Cells.Item(3,5).Value = "TC09_Canine_Filter_Breed-BorderCol_WebData.xlsx"
$ws.Cells.Item(3,2).WrapText = $true
$ws.Cells.Item(3,3).WrapText = $true
$ws.Rows.Item(3).RowHeight = 232

# --- Row 4: FilesTab ---
$ws.Cells.Item(4,1).Value = "FilesTab"
$ws.Cells.Item(4,2).Value = $filesQuery
$ws.Cells.Item(4,3).Value = $statQuery
$ws.Cells.Item(4,4).Value = "TC09_Canine_Filter_Breed-BorderCol_Neo4jData.xlsx"
$ws.Cells.Item(4,5).Value = "TC09_Canine_Filter_Breed-BorderCol_WebData.xlsx"
$ws.Cells.Item(4,2).WrapText = $true
$ws.Cells.Item(4,3).WrapText = $true
$ws.Rows.Item(4).RowHeight = 246.5

# --- Column A width + sheet view ---
$ws.Columns.Item(1).ColumnWidth = 10

[void]$ws.Range("C2").Select()
$excel.ActiveWindow.Zoom = 40

